$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 40, shifting the existing rows 40-119 down to 42-121.
$ws.Rows("40:41").Insert()

# New row 40: Macroferia Regional de Talca / Maule / Arándano (blue) / Primera
$ws.Range("A40").Value = 5
$ws.Range("B40").Value = "Macroferia Regional de Talca"
$ws.Range("C40").Value = "Maule"
$ws.Range("D40").Value = 44965
$ws.Range("E40").Value = 7
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100101
$ws.Range("H40").Value = "Berries"
$ws.Range("I40").Value = 100101001
$ws.Range("J40").Value = "Arándano (blue)"
$ws.Range("K40").Value = "Sin especificar"
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 200
$ws.Range("N40").Value = 3000
$ws.Range("O40").Value = 3000
$ws.Range("P40").Value = 3000
$ws.Range("Q40").Value = "$/bandeja 2 kilos"
$ws.Range("R40").Value = "Provincia de Curicó"
$ws.Range("S40").Value = 1500
$ws.Range("T40").Value = 2

# New row 41: Macroferia Regional de Talca / Maule / Arándano (blue) / Segunda
$ws.Range("A41").Value = 5
$ws.Range("B41").Value = "Macroferia Regional de Talca"
$ws.Range("C41").Value = "Maule"
$ws.Range("D41").Value = 44965
$ws.Range("E41").Value = 7
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100101
$ws.Range("H41").Value = "Berries"
$ws.Range("I41").Value = 100101001
$ws.Range("J41").Value = "Arándano (blue)"
$ws.Range("K41").Value = "Sin especificar"
$ws.Range("L41").Value = "Segunda"
$ws.Range("M41").Value = 150
$ws.Range("N41").Value = 2500
$ws.Range("O41").Value = 2500
$ws.Range("P41").Value = 2500
$ws.Range("Q41").Value = "$/bandeja 2 kilos"
$ws.Range("R41").Value = "Provincia de Curicó"
$ws.Range("S41").Value = 1250
$ws.Range("T41").Value = 2
